# Apply the "no-gos" docx style changes:
#  1. Add a new "Abstract Title" paragraph style (based on Normal, next
#     paragraph style Abstract).
#  2. Change the "Abstract" style's space-before from 15pt (300 twips)
#     to 5pt (100 twips).
#  3. Add a new "Footnote Block Text" paragraph style (based on the
#     built-in "Footnote Text" style).

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style -----------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$atPf = $abstractTitle.ParagraphFormat
$atPf.KeepWithNext = $true
$atPf.KeepTogether = $true
$atPf.Alignment = 1
$atPf.SpaceBefore = 15
$atPf.SpaceAfter = 0

$atFont = $abstractTitle.Font
$atFont.Size = 10
$atFont.SizeBi = 10
$atFont.Bold = $true
# wdColor is 0xBBGGRR -> target RGB 345A8A (R=0x34,G=0x5A,B=0x8A)
$atFont.Color = 9067060

# --- 2. "Abstract" style spacing tweak --------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" style -------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$fbtPf = $footnoteBlockText.ParagraphFormat
$fbtPf.SpaceBefore = 5
$fbtPf.SpaceAfter = 5
$fbtPf.FirstLineIndent = 0
$fbtPf.LeftIndent = 24
$fbtPf.RightIndent = 24

Write-Output "done"
